$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.088801741600037
$ws.Range("B1").Value = 2.455711841583252
$ws.Range("C1").Value = 5.210937976837158
$ws.Range("D1").Value = 2.215439558029175
$ws.Range("E1").Value = 1.289118885993958
